$wb = $excel.ActiveWorkbook

# Sheet: publications (column H = "Pubmed")
$wsPub = $wb.Worksheets.Item("publications")
foreach ($r in 3,6,8,12,17,20,23,24) {
    $wsPub.Range("H$r").Value = ""
}

# Sheet: preprints (column G = "Pubmed", column H = "DOI")
$wsPre = $wb.Worksheets.Item("preprints")
foreach ($r in 2,4,8,10,11,13,15,16) {
    $wsPre.Range("G$r").Value = ""
}

# Row 17 special case: G17 cleared, H17 gets the DOI value
$wsPre.Range("G17").Value = ""
$wsPre.Range("H17").Value = "10.31219/osf.io/kcvra"
